# update sim input and fix one error
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename parameter headers (row 1)
$ws.Range("B1").Value = "batch_size"
$ws.Range("C1").Value = "nboot"

# Update the sample input values (row 2) - fixes the mixed-up
# n_cases_per_iter / total_nboot values from the previous sim
$ws.Range("A2").Value = 50
$ws.Range("B2").Value = 50
$ws.Range("C2").Value = 100

# Move the active selection to D4
$ws.Range("D4").Select() | Out-Null
